$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.905.80"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").Value = "1.549.98"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.486"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("E8").Value = "  +2.96%  "

$ws.Range("E9").Value = "  -0.72%  "

$ws.Range("E10").Value = "  +0.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0854"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("D12").Value = "1.771.02"
$ws.Range("E12").Value = "  -0.27%  "

$ws.Range("D13").Value = "1.553.08"
$ws.Range("E13").Value = "  -0.15%  "

$ws.Range("E14").Value = "  +0.59%  "

$ws.Range("E15").Value = "  +0.79%  "

$ws.Range("D16").Value = "26.897.92"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.43%  "

$ws.Range("E19").Value = "  +1.42%  "

$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("E21").Value = "  -0.21%  "

$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("E24").Value = "  -0.90%  "

$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("E26").Value = "  -0.48%  "

$ws.Range("E27").Value = "  +0.38%  "

$ws.Range("E28").Value = "  +0.85%  "

$ws.Range("E29").Value = "  -0.20%  "

$ws.Range("E30").Value = "  +1.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.64%  "

$ws.Range("E32").Value = "  -0.30%  "

$ws.Range("D33").Value = "1.413.43"
$ws.Range("E33").Value = "  +3.06%  "

$ws.Range("E34").Value = "  +4.23%  "

$ws.Range("E35").Value = "  +2.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.967"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.22%  "

$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("E38").Value = "  +0.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.522"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.33%  "

$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("E41").Value = "  +4.78%  "

$ws.Range("E42").Value = "  -0.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.61%  "

$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("D47").Value = "1.684.52"
$ws.Range("E47").Value = "  -0.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0518"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.01%  "

$ws.Range("E50").Value = "  +3.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0956"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.02%  "
